$wb = $excel.ActiveWorkbook

# --- Sheet "PD" ---
$ws1 = $wb.Worksheets.Item("PD")

$ws1.Range("A2").Value = "Random Player"
$ws1.Range("B2").Value = 35
$ws1.Range("C2").Value = 4.38
$ws1.Range("D2").Value = 8

$ws1.Range("A3").Value = "Forgiving Tit for Tat"
$ws1.Range("B3").Value = 64
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 16

$ws1.Range("A4").Value = "Tit for Tat"
$ws1.Range("B4").Value = 64
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 16

$ws1.Range("A5").Value = "Grim Trigger"
$ws1.Range("B5").Value = 32
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 8

$ws1.Range("A6").Value = "Always action zero"
$ws1.Range("B6").Value = 52
$ws1.Range("C6").Value = 3.25
$ws1.Range("D6").Value = 16

# --- Sheet "PD-ASYM" ---
$ws2 = $wb.Worksheets.Item("PD-ASYM")

$ws2.Range("A2").Value = "Random Player"
$ws2.Range("B2").Value = 36
$ws2.Range("C2").Value = 4.5
$ws2.Range("D2").Value = 8

$ws2.Range("A3").Value = "Always action zero"
$ws2.Range("B3").Value = 71.59999999999999
$ws2.Range("C3").Value = 4.47
$ws2.Range("D3").Value = 16

$ws2.Range("A4").Value = "Grim Trigger"
$ws2.Range("B4").Value = 35.2
$ws2.Range("C4").Value = 4.4
$ws2.Range("D4").Value = 8

$ws2.Range("A5").Value = "Forgiving Tit for Tat"
$ws2.Range("B5").Value = 67.19999999999999
$ws2.Range("C5").Value = 4.2
$ws2.Range("D5").Value = 16

$ws2.Range("A6").Value = "Tit for Tat"
$ws2.Range("B6").Value = 67.19999999999999
$ws2.Range("C6").Value = 4.2
$ws2.Range("D6").Value = 16
